# Update column C ("Förändrad") date value from 2023-09-20 (45189) to
# 2023-09-21 (45190) for all data rows (rows 2 through 52).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}
